$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 129.2577133333333
$ws.Range("N2").Value = 387.77314
$ws.Range("O2").Value = 0.9002864985291716
$ws.Range("P2").Value = 0.9002864985291718
$ws.Range("Q2").Value = 20.66835144790445
$ws.Range("R2").Value = 186.01516303114
$ws.Range("S2").Value = 0.02383753331888343
$ws.Range("T2").Value = 0.02383753331888343

# Row 3
$ws.Range("O3").Value = 0.03306193147369147
$ws.Range("P3").Value = 0.03306193147369148
$ws.Range("R3").Value = 6.831181610799
$ws.Range("S3").Value = 0.00087540454552893
$ws.Range("T3").Value = 0.0008754045455289302

# Row 4
$ws.Range("M4").Value = 0.6895433333333333
$ws.Range("N4").Value = 2.06863
$ws.Range("O4").Value = 0.004802704126057829
$ws.Range("P4").Value = 0.00480270412605783
$ws.Range("Q4").Value = 0.1102582088477778
$ws.Range("R4").Value = 0.9923238796299999
$ws.Range("S4").Value = 0.0001271646523775262
$ws.Range("T4").Value = 0.0001271646523775263

# Row 5
$ws.Range("M5").Value = 8.231863333333333
$ws.Range("N5").Value = 24.69559
$ws.Range("O5").Value = 0.05733534367597515
$ws.Range("P5").Value = 0.05733534367597515
$ws.Range("Q5").Value = 1.316277690954444
$ws.Range("R5").Value = 11.84649921859
$ws.Range("S5").Value = 0.00151810914354327
$ws.Range("T5").Value = 0.00151810914354327

# Row 6
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6480243333333333
$ws.Range("N6").Value = 1.944073
$ws.Range("O6").Value = 0.004513522195103824
$ws.Range("P6").Value = 0.004513522195103824
$ws.Range("Q6").Value = 0.1036193069081111
$ws.Range("R6").Value = 0.9325737621729999
$ws.Range("S6").Value = 0.0001195077743441479
$ws.Range("T6").Value = 0.0001195077743441479

# Row 7
$ws.Range("M7").Value = 129.2577133333333
$ws.Range("N7").Value = 387.77314
$ws.Range("O7").Value = 0.9002864985291716
$ws.Range("P7").Value = 0.9002864985291718
$ws.Range("Q7").Value = 759.9257438590934
$ws.Range("R7").Value = 6839.33169473184
$ws.Range("S7").Value = 0.8764489652102881
$ws.Range("T7").Value = 0.8764489652102884

# Row 8
$ws.Range("O8").Value = 0.03306193147369147
$ws.Range("P8").Value = 0.03306193147369148
$ws.Range("S8").Value = 0.03218652692816254
$ws.Range("T8").Value = 0.03218652692816255

# Row 9
$ws.Range("M9").Value = 0.6895433333333333
$ws.Range("N9").Value = 2.06863
$ws.Range("O9").Value = 0.004802704126057829
$ws.Range("P9").Value = 0.00480270412605783
$ws.Range("Q9").Value = 4.053930067253333
$ws.Range("R9").Value = 36.48537060528
$ws.Range("S9").Value = 0.004675539473680302
$ws.Range("T9").Value = 0.004675539473680304

# Row 10
$ws.Range("M10").Value = 8.231863333333333
$ws.Range("N10").Value = 24.69559
$ws.Range("O10").Value = 0.05733534367597515
$ws.Range("P10").Value = 0.05733534367597515
$ws.Range("Q10").Value = 48.39637577989333
$ws.Range("R10").Value = 435.56738201904
$ws.Range("S10").Value = 0.05581723453243188
$ws.Range("T10").Value = 0.05581723453243189

# Row 11
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6480243333333333
$ws.Range("N11").Value = 1.944073
$ws.Range("O11").Value = 0.004513522195103824
$ws.Range("P11").Value = 0.004513522195103824
$ws.Range("Q11").Value = 3.809833555365334
$ws.Range("R11").Value = 34.288501998288
$ws.Range("S11").Value = 0.004394014420759675
$ws.Range("T11").Value = 0.004394014420759677
